# Updated cryptos list on Mon Apr  3 22:42:05 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# every coin row (rows 2-51) to the latest scraped snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates -------------------------------------------
# Force the range to Text format before writing so numeric-looking strings
# (e.g. "0.9997") are stored as plain text rather than being auto-coerced
# into numbers, then clear the formatting again so no number-format
# residue is left behind on the cells (they keep their original default
# style, matching the source data).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("D2").Value = "27.844.34"
$ws.Range("D3").Value = "1.806.35"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D5").Value = "306.72"
$ws.Range("D6").Value = "0.9985"
$ws.Range("D7").Value = "0.4988"
$ws.Range("D8").Value = "0.3883"
$ws.Range("D9").Value = "0.09425"
$ws.Range("D10").Value = "1.097"
$ws.Range("D11").Value = "40.43"
$ws.Range("D12").Value = "6.340"
$ws.Range("D13").Value = "0.9978"
$ws.Range("D14").Value = "20.49"
$ws.Range("D15").Value = "1.811.97"
$ws.Range("D16").Value = "7.230"
$ws.Range("D17").Value = "0.00001128"
$ws.Range("D18").Value = "92.61"
$ws.Range("D19").Value = "0.06570"
$ws.Range("D20").Value = "0.9980"
$ws.Range("D21").Value = "17.10"
$ws.Range("D22").Value = "5.935"
$ws.Range("D23").Value = "27.922.75"
$ws.Range("D24").Value = "11.05"
$ws.Range("D26").Value = "156.50"
$ws.Range("D27").Value = "20.59"
$ws.Range("D28").Value = "2.407"
$ws.Range("D29").Value = "2.010.97"
$ws.Range("D30").Value = "127.30"
$ws.Range("D31").Value = "0.1072"
$ws.Range("D32").Value = "1.053"
$ws.Range("D33").Value = "5.563"
$ws.Range("D34").Value = "3.599"
$ws.Range("D35").Value = "0.06793"
$ws.Range("D36").Value = "8.892"
$ws.Range("D37").Value = "0.02302"
$ws.Range("D38").Value = "0.2138"
$ws.Range("D40").Value = "4.926"
$ws.Range("D41").Value = "0.6207"
$ws.Range("D42").Value = "0.9979"
$ws.Range("D43").Value = "1.142"
$ws.Range("D44").Value = "13.03"
$ws.Range("D45").Value = "0.5859"
$ws.Range("D46").Value = "1.280"
$ws.Range("D47").Value = "3.668"
$ws.Range("D48").Value = "123.83"
$ws.Range("D49").Value = "1.948"
$ws.Range("D50").Value = "1.176"
$ws.Range("D2:D51").ClearFormats()

# --- Volume(1h) (column E) updates --------------------------------------
# These are already non-numeric text (leading/trailing spaces + "%"), so
# no special formatting handling is required.
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("E5").Value = "  -2.24%  "
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").Value = "  -3.98%  "
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("E9").Value = "  +20.46%  "
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("E17").Value = "  +4.41%  "
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("E28").Value = "  +3.65%  "
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("E30").Value = "  +3.88%  "
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("E35").Value = "  -6.34%  "
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("E39").Value = "  -6.39%  "
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("E46").Value = "  -6.70%  "
$ws.Range("E47").Value = "  -2.63%  "
$ws.Range("E48").Value = "  -2.98%  "
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("E50").Value = "  -4.14%  "
$ws.Range("E51").Value = "  +0.08%  "
